$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1. Add new worksheet "Sheet1" after "Favorite Colors"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet1"
$ws2.Range("C11").Value = 10
$ws2.Range("D11").Formula = "=IF(C11=5,""Yes"",""No"")"

# 2. Update the Favorite Colors source data: A5 switches from "Green" to "Red"
$ws1.Range("A5").Value = "Red"

# 3. Break a couple of the shared boolean formulas in row 5
$ws1.Range("D5").Formula = "=IF(D2<1,TRUE,FALSE)"
$ws1.Range("H5").Formula = "=H2>5"

# 4. Add new COUNTIF row for Orange at the bottom of the data
$ws1.Range("G31").Formula = "=COUNTIF(Colors,""Orange"")"

# 5. New column L width + row 13 height + the K9/K13/L13 cells
$ws1.Columns("L").ColumnWidth = 13.5
$ws1.Rows(13).RowHeight = 20

$ws1.Range("K9").Font.Color = 0

$ws1.Range("K13").Value = 10
$ws1.Range("K13").Font.Size = 16
$ws1.Range("K13").Font.Name = "Helvetica Neue"
$ws1.Range("K13").Font.Color = 1973790

$ws1.Range("L13").Formula = "=IF(K13>5,""Yes"",""No"")"
$ws1.Range("L13").Font.Size = 16
$ws1.Range("L13").Font.Name = "Helvetica Neue"
$ws1.Range("L13").Font.Color = 1973790

# 6. Remove the color-scale conditional formatting that covered C2:H2
$ws1.Range("C2:H2").FormatConditions.Delete()

# 7. Restore the original active sheet/selection
$ws1.Activate()
$ws1.Range("F5").Select()

Write-Host "done"
